# Append the new observation record as row 7 on the active ("Artfynd") sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain numeric cells.
$ws.Range("A7").Value = 111923266
$ws.Range("B7").Value = 96017
$ws.Range("E7").Value = 1001
$ws.Range("Q7").Value = 430782
$ws.Range("R7").Value = 6275341
$ws.Range("S7").Value = 10

# Plain text cells (none of these look like numbers/dates, so Excel keeps
# them as text automatically).
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "EN"
$ws.Range("F7").Value = "Flytsvalting"
$ws.Range("G7").Value = "Luronium natans"
$ws.Range("H7").Value = "(L.) Raf."
$ws.Range("J7").Value = "m²"
$ws.Range("K7").Value = "fullt utvecklade blad"
$ws.Range("N7").Value = "observerad"
$ws.Range("P7").Value = "Luroniumviken, Hängasjön, Sm"
$ws.Range("T7").Value = "Kronoberg"
$ws.Range("U7").Value = "Ljungby"
$ws.Range("V7").Value = "Småland"
$ws.Range("W7").Value = "Hamneda"
$ws.Range("AC7").Value = "Flytsvaltingarna påträffades under vatten (mycket högt vattenstånd) i ett band av 4 x 0.5 m längd."
$ws.Range("AW7").Value = "Per Ekerholm"
$ws.Range("AX7").Value = "Per Ekerholm"

# Text cells whose content would otherwise be auto-detected as a number or a
# date by Excel's normal typing heuristics - a leading single-quote forces
# them to stay plain text, matching the source data.
$ws.Range("I7").Value = "'2"
$ws.Range("Y7").Value = "'2023-09-04"
$ws.Range("AA7").Value = "'2023-09-04"

# Empty-but-present text cells.
$ws.Range("L7").Value = "'"
$ws.Range("AF7").Value = "'"
$ws.Range("AT7").Value = "'"
$ws.Range("AY7").Value = "'"

# Boolean cells.
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
